$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 290, pushing existing rows 290-352 down to 291-353.
$ws.Rows.Item(290).Insert()

$ws.Cells.Item(290, 1).Value = 5
$ws.Cells.Item(290, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(290, 3).Value = 'Maule'
$ws.Cells.Item(290, 4).Value = 44543
$ws.Cells.Item(290, 5).Value = 7
$ws.Cells.Item(290, 6).Value = 100114001
$ws.Cells.Item(290, 7).Value = 'Papa'
$ws.Cells.Item(290, 8).Value = 'Asterix'
$ws.Cells.Item(290, 9).Value = '1a nueva(o)'
$ws.Cells.Item(290, 10).Value = 1600
$ws.Cells.Item(290, 11).Value = 9000
$ws.Cells.Item(290, 12).Value = 9000
$ws.Cells.Item(290, 13).Value = 9000
$ws.Cells.Item(290, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(290, 15).Value = 'Región del Maule'
$ws.Cells.Item(290, 16).Value = 360
$ws.Cells.Item(290, 17).Value = 25
$ws.Cells.Item(290, 18).Value = 'Hortaliza'
